$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.393.74'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '2.048.75'
$ws.Range('E3').Value = '  -1.60%  '
$ws.Range('E4').Value = '  -0.14%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '229.97'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -1.84%  '
$ws.Range('E7').Value = '  +0.06%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '57.07'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -2.86%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.385'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  -2.12%  '
$ws.Range('E10').Value = '  +1.72%  '
$ws.Range('E11').Value = '  -2.03%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '14.72'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -1.93%  '
$ws.Range('D13').Value = '2.351.32'
$ws.Range('E13').Value = '  -1.58%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '20.80'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -3.04%  '
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('E16').Value = '  -1.43%  '
$ws.Range('D17').Value = '2.043.60'
$ws.Range('E17').Value = '  -1.60%  '
$ws.Range('D18').Value = '37.273.10'
$ws.Range('E18').Value = '  -1.24%  '
$ws.Range('E19').Value = '  -1.03%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '69.78'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -2.37%  '
$ws.Range('E21').Value = '  -1.11%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '226.53'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -1.69%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('E25').Value = '  -3.94%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '9.57'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -3.71%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '168.73'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -1.92%  '
$ws.Range('E28').Value = '  -5.32%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '1.38'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -2.67%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '19.00'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -2.66%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '0.118'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -2.91%  '
$ws.Range('E32').Value = '  -4.42%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '4.60'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -1.71%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '0.0615'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -3.07%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '2.43'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -1.54%  '
$ws.Range('E36').Value = '  +0.44%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('E38').Value = '  -4.69%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '5.39'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -1.32%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.0222'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -4.92%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '17.20'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +2.08%  '
$ws.Range('D42').Value = '1.498.49'
$ws.Range('E42').Value = '  +3.36%  '
$ws.Range('E43').Value = '  -1.35%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '96.64'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -5.52%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.0943'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -3.62%  '
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('E47').Value = '  -4.39%  '
$ws.Range('E48').Value = '  -4.43%  '
$ws.Range('E49').Value = '  -2.48%  '
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('D51').Value = '2.239.01'
$ws.Range('E51').Value = '  -1.53%  '
